$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 120.375
$ws.Range("I9").Value = 125.6
$ws.Range("J9").Value = 111.666664
$ws.Range("K9").Value = 125.6
$ws.Range("L9").Value = 111.666664
$ws.Range("M9").Value = 43.40000000000001
$ws.Range("N9").Value = -449.666664
$ws.Range("H12").Value = 979.55554
$ws.Range("I12").Value = 726.8333
$ws.Range("K12").Value = 726.8333
$ws.Range("M12").Value = -556.8333
$ws.Range("H18").Value = 1670.6666
$ws.Range("I18").Value = 1758.909
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 1758.909
$ws.Range("L18").Value = 700
$ws.Range("M18").Value = -1474.909
$ws.Range("N18").Value = -1268
$ws.Range("H40").Value = 2061.5386
$ws.Range("J40").Value = 2144.4443
$ws.Range("L40").Value = 2144.4443
$ws.Range("N40").Value = -2494.4443
$ws.Range("H70").Value = 368486
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 368486
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 1105458
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -1105998
$ws.Range("H73").Value = 368486
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 368486
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 1105458
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -1107330
$ws.Range("H107").Value = 607
$ws.Range("I107").Value = 166
$ws.Range("K107").Value = 166
$ws.Range("M107").Value = 1754
$ws.Range("H116").Value = 5834.778
$ws.Range("J116").Value = 6831.375
$ws.Range("L116").Value = 6831.375
$ws.Range("N116").Value = -13715.375
$ws.Range("H129").Value = 52098.5
$ws.Range("I129").Value = 1197
$ws.Range("J129").Value = 103000
$ws.Range("K129").Value = 3591
$ws.Range("L129").Value = 309000
$ws.Range("M129").Value = 1409
$ws.Range("N129").Value = -319000
$ws.Range("H132").Value = 1255.3448
$ws.Range("I132").Value = 1286.1072
$ws.Range("K132").Value = 3858.3216
$ws.Range("M132").Value = -1328.3216
$ws.Range("H137").Value = 2925.9443
$ws.Range("J137").Value = 5099.8
$ws.Range("L137").Value = 15299.4
$ws.Range("N137").Value = -20399.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1838.75
$ws.Range("I45").Value = 1838.75
$ws.Range("K45").Value = 1838.75
$ws.Range("M45").Value = -1461.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4403.2
$ws.Range("I105").Value = 3170.2778
$ws.Range("K105").Value = 3170.2778
$ws.Range("M105").Value = -1423.2778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 946.5
$ws.Range("I16").Value = 949.75
$ws.Range("K16").Value = 949.75
$ws.Range("M16").Value = -662.75
$ws.Range("H22").Value = 547.4167
$ws.Range("I22").Value = 483.5
$ws.Range("J22").Value = 675.25
$ws.Range("K22").Value = 483.5
$ws.Range("L22").Value = 675.25
$ws.Range("M22").Value = -133.5
$ws.Range("N22").Value = -1375.25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = $null
$ws.Range("H31").Value = 5823.857
$ws.Range("I31").Value = 2499.25
$ws.Range("J31").Value = 7153.7
$ws.Range("K31").Value = 2499.25
$ws.Range("L31").Value = 7153.7
$ws.Range("M31").Value = -2204.25
$ws.Range("N31").Value = -7743.7
$ws.Range("H34").Value = 5823.857
$ws.Range("I34").Value = 2499.25
$ws.Range("J34").Value = 7153.7
$ws.Range("K34").Value = 2499.25
$ws.Range("L34").Value = 7153.7
$ws.Range("M34").Value = -2297.25
$ws.Range("N34").Value = -7557.7
$ws.Range("H62").Value = 59713.43
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 82398.8
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 82398.8
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -83646.8
$ws.Range("H65").Value = 59713.43
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 82398.8
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 411994
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -418234
$ws.Range("H99").Value = 17256.354
$ws.Range("I99").Value = 15373
$ws.Range("K99").Value = 15373
$ws.Range("M99").Value = -13875
$ws.Range("H105").Value = 1286
$ws.Range("I105").Value = 766.5
$ws.Range("J105").Value = 2325
$ws.Range("K105").Value = 766.5
$ws.Range("L105").Value = 2325
$ws.Range("M105").Value = 980.5
$ws.Range("N105").Value = -5819
$ws.Range("H107").Value = 650.13635
$ws.Range("J107").Value = 1037.2858
$ws.Range("L107").Value = 1037.2858
$ws.Range("N107").Value = -4877.2858
$ws.Range("H113").Value = 946.5
$ws.Range("I113").Value = 949.75
$ws.Range("K113").Value = 949.75
$ws.Range("M113").Value = 1220.25
$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178
$ws.Range("H122").Value = 4409.875
$ws.Range("I122").Value = 4507.3335
$ws.Range("J122").Value = 4117.5
$ws.Range("K122").Value = 13522.0005
$ws.Range("L122").Value = 12352.5
$ws.Range("M122").Value = -11072.0005
$ws.Range("N122").Value = -17252.5
$ws.Range("H126").Value = 17256.354
$ws.Range("I126").Value = 15373
$ws.Range("K126").Value = 46119
$ws.Range("M126").Value = -43649
$ws.Range("H134").Value = 3415.7058
$ws.Range("I134").Value = 2727.818
$ws.Range("J134").Value = 4676.8335
$ws.Range("K134").Value = 8183.454000000001
$ws.Range("L134").Value = 14030.5005
$ws.Range("M134").Value = -5648.454000000001
$ws.Range("N134").Value = -19100.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66770.39999999999
$ws.Range("I2").Value = 90993.91
$ws.Range("J2").Value = 155.75
$ws.Range("K2").Value = 545963.46
$ws.Range("L2").Value = 934.5
$ws.Range("M2").Value = -545850.46
$ws.Range("N2").Value = -1160.5
$ws.Range("H23").Value = 343.33334
$ws.Range("I23").Value = 375
$ws.Range("K23").Value = 1125
$ws.Range("M23").Value = -890
$ws.Range("H80").Value = 4446.1665
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 5992.3335
$ws.Range("K80").Value = 8700
$ws.Range("L80").Value = 17977.0005
$ws.Range("M80").Value = -7764
$ws.Range("N80").Value = -19849.0005
$ws.Range("H83").Value = 4446.1665
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 5992.3335
$ws.Range("K83").Value = 26100
$ws.Range("L83").Value = 53931.0015
$ws.Range("M83").Value = -21420
$ws.Range("N83").Value = -63291.0015
$ws.Range("H86").Value = 171.33333
$ws.Range("I86").Value = 113.333336
$ws.Range("K86").Value = 340.000008
$ws.Range("M86").Value = 845.999992
$ws.Range("H89").Value = 171.33333
$ws.Range("I89").Value = 113.333336
$ws.Range("K89").Value = 1020.000024
$ws.Range("M89").Value = 4907.999976
$ws.Range("H92").Value = 499
$ws.Range("I92").Value = 499
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1497
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -249
$ws.Range("N92").Value = $null
$ws.Range("H97").Value = 141.5
$ws.Range("I97").Value = 141.5
$ws.Range("K97").Value = 424.5
$ws.Range("M97").Value = 71.5
$ws.Range("H98").Value = 775
$ws.Range("I98").Value = 775
$ws.Range("K98").Value = 2325
$ws.Range("M98").Value = -827
$ws.Range("H107").Value = 623.0278
$ws.Range("J107").Value = 603.3143
$ws.Range("L107").Value = 1809.9429
$ws.Range("N107").Value = -5649.9429
$ws.Range("H113").Value = 1217.7142
$ws.Range("J113").Value = 1094.4445
$ws.Range("L113").Value = 3283.3335
$ws.Range("N113").Value = -7623.333500000001
$ws.Range("H117").Value = 600
$ws.Range("I117").Value = 633.3333
$ws.Range("J117").Value = 500
$ws.Range("K117").Value = 1899.9999
$ws.Range("L117").Value = 1500
$ws.Range("M117").Value = 1542.0001
$ws.Range("N117").Value = -8384
$ws.Range("H122").Value = 362.7619
$ws.Range("I122").Value = 312
$ws.Range("K122").Value = 2808
$ws.Range("M122").Value = -358
$ws.Range("H137").Value = 7015.4443
$ws.Range("I137").Value = 4999.5
$ws.Range("J137").Value = 8023.4165
$ws.Range("K137").Value = 14998.5
$ws.Range("L137").Value = 24070.2495
$ws.Range("M137").Value = -9898.5
$ws.Range("N137").Value = -34270.24950000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4678.1665
$ws.Range("J61").Value = 3942.889
$ws.Range("L61").Value = 3942.889
$ws.Range("N61").Value = -4346.889
$ws.Range("H93").Value = 289.83334
$ws.Range("I93").Value = 291.8
$ws.Range("K93").Value = 291.8
$ws.Range("M93").Value = 956.2
$ws.Range("H113").Value = 4678.1665
$ws.Range("J113").Value = 3942.889
$ws.Range("L113").Value = 3942.889
$ws.Range("N113").Value = -8282.888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7406.3335
$ws.Range("I62").Value = 5250
$ws.Range("J62").Value = 8484.5
$ws.Range("K62").Value = 5250
$ws.Range("L62").Value = 8484.5
$ws.Range("M62").Value = -4626
$ws.Range("N62").Value = -9732.5
$ws.Range("H65").Value = 7406.3335
$ws.Range("I65").Value = 5250
$ws.Range("J65").Value = 8484.5
$ws.Range("K65").Value = 26250
$ws.Range("L65").Value = 42422.5
$ws.Range("M65").Value = -23130
$ws.Range("N65").Value = -48662.5
$ws.Range("H107").Value = 865.2857
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 865.2857
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2595.8571
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = -6435.8571
$ws.Range("H136").Value = 2973.8696
$ws.Range("I136").Value = 823.5833
$ws.Range("J136").Value = 5319.636
$ws.Range("K136").Value = 2470.7499
$ws.Range("L136").Value = 15958.908
$ws.Range("M136").Value = 79.2501000000002
$ws.Range("N136").Value = -21058.908
